$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Helper: write a numeric-looking string into a cell and have it survive
# as literal text (matching the workbook's convention of storing the
# fund-metric columns as text) instead of Excel's default behaviour of
# silently coercing a numeric-looking string typed into a General-format
# cell into a real number. Writing through a text-formatted scratch cell
# and pasting *values only* keeps the destination cell's own style/format
# untouched (General, no quote-prefix marker).
# -----------------------------------------------------------------------
function Set-TextValue($ws, $range, [string]$text) {
    $helper = $ws.Range("ZZ1")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $helper.Clear()
    $excel.CutCopyMode = $false
}

# -----------------------------------------------------------------------
# 1) "总计" (Total) sheet: insert a new row for the 2022-Q3 quarter right
#    after the header, shifting the existing quarter rows down by one.
# -----------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# the freshly inserted row inherits a blended style from Excel's own
# auto-fill-format heuristics; re-copy the plain formatting from the row
# below (still exactly as authored) so column A keeps its centred/bordered
# style and columns B:D stay unstyled, like every other data row.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$totalData = @(
    @(0, "2022-Q3", 3, 0.05),
    @(1, "2022-Q2", 3, 0.04),
    @(2, "2022-Q1", 5, 0.19),
    @(3, "2021-Q4", 5, 1.28),
    @(4, "2021-Q3", 4, 0.02)
)

$r = 2
foreach ($row in $totalData) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# -----------------------------------------------------------------------
# 2) Add a brand-new "2022-Q3" sheet, positioned right after "总计" (i.e.
#    right before the existing "2022-Q2" tab). It holds the same fund
#    line-up as "2022-Q2" but with refreshed position metrics.
# -----------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

$q3Data = @(
    @("0.64", "80.13", "2.49", "0.0159"),
    @("0.64", "80.13", "2.49", "0.0159"),
    @("0.64", "80.13", "2.49", "0.0159")
)

$r = 2
foreach ($row in $q3Data) {
    Set-TextValue $q3 $q3.Cells.Item($r, 4) $row[0]
    Set-TextValue $q3 $q3.Cells.Item($r, 5) $row[1]
    Set-TextValue $q3 $q3.Cells.Item($r, 6) $row[2]
    Set-TextValue $q3 $q3.Cells.Item($r, 7) $row[3]
    $r++
}

Write-Host "edit applied"
